$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on sheet "Hoja1" (cell A1) ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$check = [char]0x2705

$oldLine1 = "$check 1000 Bs = 2.33 = 8796.84 pesos"
$newLine1 = "$check 1000 Bs = 2.29 = 8675.28 pesos"
$oldLine2 = "$check 8796.84 pesos = 2.33 = 969.09 Bs"
$newLine2 = "$check 8675.28 pesos = 2.27 = 955.6 Bs"

[string]$text = $ws1.Range("A1").Value2
$text = $text.Replace($oldLine1, $newLine1)
$text = $text.Replace($oldLine2, $newLine2)
$ws1.Range("A1").Value = $text

# --- Update the rate figures on sheet "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 436.99
$ws2.Range("O10").Value = 3791.01
$ws2.Range("N12").Value = 3822
$ws2.Range("O12").Value = 421
